# Fixed objectives on every mission.
#
# Adds a new localization entry (key / English default / Portuguese
# translation) for the "light up torches" objective, appended right after
# the last existing data row on the "Localization" sheet:
#   A72 = Objective.LightTorches
#   B72 = Light up torches
#   C72 = Acenda as tochas

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# Insert a new row at 72, shifting nothing below it (there is nothing below)
# but — crucially — inheriting the formatting of the row above (row 71),
# exactly like Excel does when a row is inserted via the UI. This is what
# gives B72/C72 the same cell style as the rest of the table (applyFont
# cell style used throughout column B/C) while leaving column A unstyled,
# matching the rest of the sheet.
$ws.Rows.Item(72).Insert(-4121, 0)

$ws.Cells.Item(72, 1).Value = "Objective.LightTorches"
$ws.Cells.Item(72, 2).Value = "Light up torches"
$ws.Cells.Item(72, 3).Value = "Acenda as tochas"

# A stray formatted (but empty) cell was left behind in C73 — same kind of
# leftover formatting artifact already present elsewhere on this sheet
# (e.g. D22). Recreate it with the underlined font style.
$ws.Range("C73").Font.Underline = 2

# Restore the recorded selection / scroll position from the edit session.
$ws.Range("A72").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 43 | Out-Null
